# Comp230/merged_sales_data.xlsx - remove duplicate rows that were
# introduced while merging the two source spreadsheets, and clean up a
# couple of truly-blank rows on the "invoices" sheet.
#
# customers common to both source spreadsheets only show once in the
# merged target; next step is removing the (now) blank rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "customers": rows 18-25 duplicate customer records 123-130
# that already exist in rows 2,3,5,6,8,9,10,11. Clear their contents
# (NOT a row-delete/shift, the customers below at rows 27-30 stay put).
# ---------------------------------------------------------------
$wsCustomers = $wb.Worksheets.Item("customers")
$wsCustomers.Range("A18:C25").ClearContents()

# ---------------------------------------------------------------
# Sheet "invoices": two genuinely empty rows (55 and 56) get removed
# with the standard shift-up delete, pulling the rest of the invoice
# rows (old 58..103) up to new rows 56..101.
# ---------------------------------------------------------------
$wsInvoices = $wb.Worksheets.Item("invoices")
$wsInvoices.Rows("55:56").Delete()

# ---------------------------------------------------------------
# Sheet "products": rows 11-18 duplicate product records 1-8 that
# already exist in rows 2-9. Clear their contents (rows below don't
# exist, so this just shrinks the used range down to row 9).
# ---------------------------------------------------------------
$wsProducts = $wb.Worksheets.Item("products")
$wsProducts.Range("A11:D18").ClearContents()
